$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "62.942.47"
$ws.Range("E2").Value = "  -5.73%  "
$ws.Range("D3").Value = "3.259.98"
$ws.Range("E3").Value = "  -6.79%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "179.98"
$ws.Range("E5").Value = "  -10.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "518.31"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "3.261.20"
$ws.Range("E8").Value = "  -6.55%  "
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.71"
$ws.Range("E11").Value = "  -4.54%  "
$ws.Range("E12").Value = "  -8.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  -6.07%  "
$ws.Range("E14").Value = "  -6.80%  "
$ws.Range("D15").Value = "3.777.03"
$ws.Range("E15").Value = "  -7.30%  "
$ws.Range("E16").Value = "  -5.62%  "
$ws.Range("D17").Value = "3.263.74"
$ws.Range("E17").Value = "  -6.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.53"
$ws.Range("E18").Value = "  -4.78%  "
$ws.Range("D19").Value = "62.916.01"
$ws.Range("E19").Value = "  -5.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.92"
$ws.Range("E20").Value = "  -7.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.948"
$ws.Range("E21").Value = "  -7.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "371.04"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("E23").Value = "  -5.31%  "
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.76"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.83"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.07"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.41"
$ws.Range("E28").Value = "  -3.98%  "
$ws.Range("E29").Value = "  -5.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.31"
$ws.Range("E30").Value = "  -5.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.43"
$ws.Range("E31").Value = "  -6.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.80"
$ws.Range("E32").Value = "  -6.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "634.60"
$ws.Range("E33").Value = "  -5.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.27"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.65"
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.400"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.33"
$ws.Range("E39").Value = "  -7.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "2.959.40"
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.125"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "0.0₃0653"
$ws.Range("E43").Value = "  -6.80%  "
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("E45").Value = "  -12.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0393"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.58"
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.80"
$ws.Range("E48").Value = "  +8.00%  "
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.98"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("E51").Value = "  -11.83%  "
